$wb = $excel.ActiveWorkbook

# --- Add the new "checkpoint" sheet right after Sheet1 ---
$sheet1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "checkpoint"

# --- Fill in the header / checkpoint data (two identical rows) ---
$headers = @("CURA Healthcare Service", "Make Appointment", "Appointment Confirmation", "History", "CURA Healthcare Service")
for ($col = 1; $col -le 5; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headers[$col - 1]
    $ws2.Cells.Item(2, $col).Value = $headers[$col - 1]
}

# --- Wrap text on the columns that hold the longer labels (not column C) ---
$ws2.Range("A1:B2").WrapText = $true
$ws2.Range("D1:E2").WrapText = $true

# --- Column widths roughly matching the recorded layout ---
$ws2.Columns.Item(1).ColumnWidth = 22.25
$ws2.Columns.Item(2).ColumnWidth = 22.09
$ws2.Columns.Item(3).ColumnWidth = 22.09
$ws2.Columns.Item(5).ColumnWidth = 22.75

# --- Page setup mirrors the rest of the workbook (same margins/header/footer) ---
$ws2.PageSetup.LeftMargin = 56.7
$ws2.PageSetup.RightMargin = 56.7
$ws2.PageSetup.TopMargin = 75.8
$ws2.PageSetup.BottomMargin = 75.8
$ws2.PageSetup.HeaderMargin = 56.7
$ws2.PageSetup.FooterMargin = 56.7
$ws2.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ws2.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12Page &P"
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selection/view state on the new sheet ---
$ws2.Range("F2:F3").Select() | Out-Null

# --- Sheet1 keeps its own B2 selection, plus a lingering F2:F3 area ---
$u = $excel.Union($sheet1.Range("F2:F3"), $sheet1.Range("B2"))
$u.Select() | Out-Null
$sheet1.Range("B2").Select() | Out-Null

# checkpoint is the active tab once everything is in place
$ws2.Activate() | Out-Null
